# B6-PowerPoint.pptx — update the table style applied to the three data
# tables (slides 14, 15 and 16) from the old tableStyleId
# {A62A2CE7-E7AF-4133-A4F9-6FD590943077} to the new one
# {BACCF9CB-AD5A-41AD-96EC-084BE2D7BD3A}, matching the Design > Table
# Styles gallery choice made by the author.

$OldStyleId = "{A62A2CE7-E7AF-4133-A4F9-6FD590943077}"
$NewStyleId = "{BACCF9CB-AD5A-41AD-96EC-084BE2D7BD3A}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if ($shape.HasTable) {
            $table = $shape.Table

            if ($table.Style.Name -eq $OldStyleId) {
                $table.ApplyStyle($NewStyleId)
            }
        }
    }
}
